# Apply the "Add files via upload" edit to mamografias_CA_mama.xlsx.
#
# What changed (per the OOXML diff):
#   1. On sheet "cancer_de_mama", a batch of cells that only held throwaway
#      placeholder text ("-" or "_") are cleared back to truly empty cells.
#      Clearing them is also what causes those two now-unused shared
#      strings to drop out of sharedStrings.xml on save (and every other
#      shared-string index used on the sheet to shift down accordingly) -
#      that part is automatic bookkeeping, not something to hand-edit.
#   2. The sheet's view no longer scrolls to keep column AS in view
#      (topLeftCell="AS1" is gone) and the selection moves from the single
#      cell A11 to the full row range A2:XFD12.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("cancer_de_mama")

# 1. Clear the placeholder cells (previously "-" or "_") back to empty.
$cellsToClear = "I2","J2","J3","I4","J4","I7","J7","I8","J8","I9","J9","J11","BK12","BL12","BM12","BN12","BO12","BP12","BQ12","BR12"
foreach ($addr in $cellsToClear) {
    $ws.Range($addr).ClearContents()
}

# 2. Update the sheet view: reset horizontal scroll (drop topLeftCell) and
#    select A2:XFD12 (whole rows 2-12) with A2 as the active cell.
$ws.Activate()
$ws.Range("A2:XFD12").Select()
